$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed rows 47 and 48 by copying the formatting (and placeholder values) of row 46 ---
$ws.Range("A46:P46").Copy($ws.Range("A47:P47"))
$ws.Range("A46:P46").Copy($ws.Range("A48:P48"))
$ws.Rows.Item(47).RowHeight = 17
$ws.Rows.Item(48).RowHeight = 17

# Row 46's "H" column carries the vertical-centered style (s=6); the two new rows use
# the plain style (s=5) for column H instead, so pull the format from a neighboring
# s=5 cell.
$ws.Range("I47").Copy()
$ws.Range("H47").PasteSpecial(-4122)
$ws.Range("I48").Copy()
$ws.Range("H48").PasteSpecial(-4122)

# Row 46 has content in column L ("11--12"); the new rows must not carry that over.
$ws.Range("L47").ClearContents()
$ws.Range("L48").ClearContents()

# --- Populate the two new (joke) course rows ---
$ws.Cells.Item(47, 2).Value = "我只是来求课的TvT"
$ws.Cells.Item(48, 2).Value = "我是雷锋我就是来出课的"
$ws.Cells.Item(47, 15).Value = "理教666"
$ws.Cells.Item(47, 8).Value = "1--12"
$ws.Cells.Item(47, 7).Value = "无"

$ws.Cells.Item(47, 1).Value = 31
$ws.Cells.Item(47, 3).Value = "选修"
$ws.Cells.Item(47, 4).Value = "原理"
$ws.Cells.Item(47, 5).Value = 2
$ws.Cells.Item(47, 6).Value = "1--16"

$ws.Cells.Item(48, 1).Value = 32
$ws.Cells.Item(48, 3).Value = "选修"
$ws.Cells.Item(48, 4).Value = "原理"
$ws.Cells.Item(48, 5).Value = 2
$ws.Cells.Item(48, 6).Value = "1--16"
$ws.Cells.Item(48, 7).Value = "无"
$ws.Cells.Item(48, 8).Value = "1--12"
$ws.Cells.Item(48, 15).Value = "理教666"

# Columns I, J, K, M, N, P stay empty for both new rows (already blank after the copy).

# --- Update the sheet selection to match the edited state ---
$ws.Range("J50").Select()

# --- Page setup: the saved file now also carries explicit page setup info ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
